# Add an "SW Efficiency" column to the Pluvial_Domain sheet, between the
# existing "SW Capacity (in)" column and the "Comments" column, and
# populate the two rows that have SW Infrastructure = YES with a value of 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pluvial_Domain")

# Insert a new column at G (shifts the old "Comments" column from G to H).
$ws.Columns.Item(7).Insert()

# Header for the new column.
$ws.Range("G1").Value = "SW Efficiency"

# Values for the rows that previously had SW data (rows 2 and 4).
$ws.Range("G2").Value = 1
$ws.Range("G4").Value = 1

# Give the new column the same display width as the neighbouring
# "SW Capacity (in)" column (F).
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# Update the view: Pluvial_Domain becomes the active/selected sheet,
# with the new G4 cell selected.
$ws.Activate()
$ws.Range("G4").Select()
